$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '29.236.44'
Set-TextValue $ws.Range('E2') '  +0.39%  '
Set-TextValue $ws.Range('D3') '1.858.39'
Set-TextValue $ws.Range('E3') '  +0.35%  '
Set-TextValue $ws.Range('D4') '0.9997'
Set-TextValue $ws.Range('E4') '  -0.05%  '
Set-TextValue $ws.Range('D5') '0.7106'
Set-TextValue $ws.Range('E5') '  +2.09%  '
Set-TextValue $ws.Range('D6') '238.23'
Set-TextValue $ws.Range('E6') '  -0.10%  '
Set-TextValue $ws.Range('D7') '0.9997'
Set-TextValue $ws.Range('E7') '  -0.05%  '
Set-TextValue $ws.Range('D8') '0.08008'
Set-TextValue $ws.Range('E8') '  +4.73%  '
Set-TextValue $ws.Range('D9') '0.3035'
Set-TextValue $ws.Range('E9') '  +0.19%  '
Set-TextValue $ws.Range('D10') '23.54'
Set-TextValue $ws.Range('E10') '  +0.73%  '
Set-TextValue $ws.Range('D11') '0.08193'
Set-TextValue $ws.Range('E11') '  +0.84%  '
Set-TextValue $ws.Range('D12') '1.835.70'
Set-TextValue $ws.Range('E12') '  +0.99%  '
Set-TextValue $ws.Range('D13') '5.175'
Set-TextValue $ws.Range('D14') '0.7043'
Set-TextValue $ws.Range('E14') '  -3.07%  '
Set-TextValue $ws.Range('D15') '89.71'
Set-TextValue $ws.Range('E15') '  +0.69%  '
Set-TextValue $ws.Range('D16') '29.238.15'
Set-TextValue $ws.Range('E16') '  +0.44%  '
Set-TextValue $ws.Range('D17') '5.830'
Set-TextValue $ws.Range('D18') '0.000007872'
Set-TextValue $ws.Range('E18') '  +1.90%  '
Set-TextValue $ws.Range('D19') '13.27'
Set-TextValue $ws.Range('E19') '  +0.44%  '
Set-TextValue $ws.Range('D20') '238.00'
Set-TextValue $ws.Range('E20') '  +0.45%  '
Set-TextValue $ws.Range('D21') '0.9995'
Set-TextValue $ws.Range('E21') '  +0.01%  '
Set-TextValue $ws.Range('D22') '2.097.45'
Set-TextValue $ws.Range('E22') '  +0.17%  '
Set-TextValue $ws.Range('D23') '1.000'
Set-TextValue $ws.Range('E23') '  -0.03%  '
Set-TextValue $ws.Range('D24') '7.443'
Set-TextValue $ws.Range('E24') '  -2.23%  '
Set-TextValue $ws.Range('D25') '162.73'
Set-TextValue $ws.Range('E25') '  +0.95%  '
Set-TextValue $ws.Range('E26') '  -0.55%  '
Set-TextValue $ws.Range('D27') '0.1448'
Set-TextValue $ws.Range('E27') '  +0.34%  '
Set-TextValue $ws.Range('D28') '18.11'
Set-TextValue $ws.Range('E28') '  +0.35%  '
Set-TextValue $ws.Range('D29') '1.926'
Set-TextValue $ws.Range('E29') '  -3.20%  '
Set-TextValue $ws.Range('D30') '1.432'
Set-TextValue $ws.Range('E30') '  +1.62%  '
Set-TextValue $ws.Range('E31') '  -0.51%  '
Set-TextValue $ws.Range('D32') '4.367'
Set-TextValue $ws.Range('E32') '  -2.66%  '
Set-TextValue $ws.Range('D33') '4.014'
Set-TextValue $ws.Range('E34') '  -0.33%  '
Set-TextValue $ws.Range('D35') '1.161'
Set-TextValue $ws.Range('E35') '  -2.33%  '
Set-TextValue $ws.Range('D36') '0.7101'
Set-TextValue $ws.Range('E36') '  +1.26%  '
Set-TextValue $ws.Range('D37') '0.9975'
Set-TextValue $ws.Range('E37') '  -2.67%  '
Set-TextValue $ws.Range('E38') '  +0.74%  '
Set-TextValue $ws.Range('D39') '0.01855'
Set-TextValue $ws.Range('E39') '  +0.01%  '
Set-TextValue $ws.Range('D40') '2.716'
Set-TextValue $ws.Range('E40') '  +1.44%  '
Set-TextValue $ws.Range('D41') '0.9295'
Set-TextValue $ws.Range('E41') '  -0.37%  '
Set-TextValue $ws.Range('D42') '1.130.58'
Set-TextValue $ws.Range('E42') '  +4.74%  '
Set-TextValue $ws.Range('D43') '0.4258'
Set-TextValue $ws.Range('E43') '  -0.16%  '
Set-TextValue $ws.Range('B44') 'FraxShare'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D44') '5.866'
Set-TextValue $ws.Range('E44') '  -1.72%  '
Set-TextValue $ws.Range('B45') 'Aave'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '70.45'
Set-TextValue $ws.Range('E45') '  +0.28%  '
Set-TextValue $ws.Range('D46') '0.9993'
Set-TextValue $ws.Range('E46') '  -0.10%  '
Set-TextValue $ws.Range('D47') '102.86'
Set-TextValue $ws.Range('E47') '  -0.21%  '
Set-TextValue $ws.Range('D48') '0.5342'
Set-TextValue $ws.Range('E48') '  -4.24%  '
Set-TextValue $ws.Range('E49') '  -0.63%  '
Set-TextValue $ws.Range('B50') 'RocketPoolETH'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D50') '1.989.89'
Set-TextValue $ws.Range('E50') '  +0.04%  '
Set-TextValue $ws.Range('B51') 'EnergySwap'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '9.164'
Set-TextValue $ws.Range('E51') '  -0.22%  '
